# Vincule tabla_contactos a ventanaPrincipal y agregue librerias en readme
#
# Sync the sample "tabla_contactos" sheet with the data now bound to
# ventanaPrincipal: refresh Sonia Hernandez's row, drop the old "cvc/vv"
# scratch row, replace the "few/ff" scratch row with her real contact
# card, add Graciela Gomez, and leave a fresh blank row ready for the
# next contact entered from the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Hernandez, Sonia -> refreshed phone/email/address ---------
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "3644123456"
$ws.Range("E2").Value = "srhernandez@gmail.com"
$ws.Range("F2").Value = "Bn 234"

# --- Row 3: drop the old "cvc/vv" scratch/test contact -----------------
$ws.Range("A3:F3").ClearContents()

# --- Row 4: replace the "few/ff" scratch contact with the real card ---
$ws.Range("A4").Value = "Hernandez"
$ws.Range("B4").Value = "Sonia"
$ws.Range("C4").Value = "Soni"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "3644567800"
$ws.Range("E4").Value = "shernandez@gmail.com"
$ws.Range("F4").Value = "Ant Arg"

# --- Row 6: new contact, Gomez, Graciela --------------------------------
$ws.Range("A6").Value = "Gomez"
$ws.Range("B6").Value = "Graciela"
$ws.Range("C6").Value = "Gra"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "34567890"
$ws.Range("E6").Value = "gra12@gmail.com"
$ws.Range("F6").Value = "La Rioja"

# --- Row 7: leave a blank row ready for the next entry ------------------
$ws.Range("A7:F7").NumberFormat = "@"
$ws.Range("A7:F7").Value = ""
